# Update "Countries & provincias Spain" data.
# - Refreshes the "Datos actualizados..." timestamp banner.
# - Re-sorts/updates the province case counts (columns B:E) with the
#   latest figures, which also shuffles several province names up/down
#   in column A as their totals overtake neighbours, and introduces a
#   few provinces that newly cross into the visible ranking
#   (Illes Balears, Pontevedra, Ourense, Cantabria, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner above the table.
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 13:16"

# Row, Provincia, Casos totales, Casos activos, Recuperados, Muertes
$data = @(
    @(4, 'Madrid', 7165, 1186, 5351, 628),
    @(5, 'Cataluña', 3270, 3, 3185, 82),
    @(6, 'Araba/Alava', 703, 21, 655, 48),
    @(7, 'Valencia/Valencia', 627, 7, 608, 12),
    @(8, 'Navarra', 554, 2, 546, 6),
    @(9, 'Bizkaia/Vizcaya', 539, 21, 522, 17),
    @(10, 'La Rioja', 497, 13, 477, 7),
    @(11, 'Alacant/Alicante', 372, 11, 350, 11),
    @(12, 'Malaga', 361, 0, 348, 13),
    @(13, 'Asturias', 344, 12, 329, 3),
    @(14, 'Toledo', 293, 14, 265, 14),
    @(15, 'A Coruña', 270, 0, 267, 3),
    @(16, 'Albacete', 259, 8, 231, 20),
    @(17, 'Zaragoza', 224, 0, 210, 14),
    @(18, 'Gipuzkoa/Guipuzcoa', 223, 21, 217, 6),
    @(19, 'Ciudad Real', 216, 6, 189, 21),
    @(20, 'Murcia', 215, 1, 214, 0),
    @(21, 'Tenerife', 210, 6, 135, 3),
    @(22, 'Caceres', 206, 2, 194, 10),
    @(23, 'Guadalajara', 205, 2, 200, 3),
    @(24, 'Illes Balears', 203, 10, 189, 4),
    @(25, 'Pontevedra', 193, 0, 191, 2),
    @(26, 'Burgos', 187, 14, 165, 8),
    @(27, 'Granada', 176, 0, 169, 7),
    @(28, 'Aragon', 174, 0, 163, 11),
    @(29, 'Illes Balears*', 169, 6, 161, 2),
    @(30, 'Salamanca', 149, 8, 129, 12),
    @(31, 'Cantabria', 144, 10, 133, 1),
    @(32, 'Leon', 134, 3, 126, 5),
    @(33, 'Sevilla', 133, 1, 131, 1),
    @(34, 'Segovia', 121, 3, 111, 7),
    @(35, 'Valladolid', 115, 1, 111, 3),
    @(36, 'Castello/Castellon', 104, 1, 102, 1),
    @(37, 'Cordoba', 101, 0, 101, 0),
    @(38, 'Badajoz', 91, 5, 86, 0),
    @(39, 'Jaen', 87, 0, 85, 2),
    @(40, 'Cadiz', 84, 0, 84, 0),
    @(41, 'Cuenca', 72, 4, 64, 4),
    @(42, 'Ourense', 63, 0, 63, 0),
    @(43, 'Avila', 59, 2, 55, 2),
    @(44, 'Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena', 58, 0, 58, 3),
    @(45, 'Soria', 58, 4, 49, 5),
    @(46, 'Gran Canaria', 55, 0, 54, 1),
    @(47, 'Lugo', 46, 0, 45, 1),
    @(48, 'Almeria', 37, 0, 37, 0),
    @(49, 'Zamora', 31, 1, 29, 1),
    @(50, 'Teruel', 27, 0, 26, 1),
    @(51, 'Huesca', 24, 0, 24, 0),
    @(52, 'Melilla', 24, 0, 24, 0),
    @(53, 'Huelva', 23, 0, 23, 0),
    @(54, 'Palencia', 14, 1, 13, 0),
    @(55, 'Fuerteventura', 11, 0, 11, 0),
    @(56, 'Arroyo de la Luz', 7, 0, 7, 0),
    @(57, 'Ceuta', 5, 0, 5, 0),
    @(58, 'La Palma', 5, 0, 5, 0),
    @(59, 'Lanzarote', 3, 0, 3, 0),
    @(60, 'La Gomera', 3, 2, 1, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}
